$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "29.170.54"
Set-TextValue $ws.Range("E2") "  -0.44%  "

Set-TextValue $ws.Range("D3") "1.825.87"
Set-TextValue $ws.Range("E3") "  -0.73%  "

Set-TextValue $ws.Range("E4") "  +0.41%  "

Set-TextValue $ws.Range("D5") "234.80"
Set-TextValue $ws.Range("E5") "  -1.70%  "

Set-TextValue $ws.Range("D6") "0.5976"
Set-TextValue $ws.Range("E6") "  -4.35%  "

Set-TextValue $ws.Range("E7") "  +0.31%  "

Set-TextValue $ws.Range("D8") "0.06929"
Set-TextValue $ws.Range("E8") "  -5.97%  "

Set-TextValue $ws.Range("D9") "0.2740"
Set-TextValue $ws.Range("E9") "  -5.11%  "

Set-TextValue $ws.Range("D10") "23.27"
Set-TextValue $ws.Range("E10") "  -6.01%  "

Set-TextValue $ws.Range("D11") "0.07600"
Set-TextValue $ws.Range("E11") "  -1.49%  "

Set-TextValue $ws.Range("D12") "1.833.02"
Set-TextValue $ws.Range("E12") "  +0.07%  "

Set-TextValue $ws.Range("D13") "4.726"
Set-TextValue $ws.Range("E13") "  -4.46%  "

Set-TextValue $ws.Range("D14") "0.6234"
Set-TextValue $ws.Range("E14") "  -5.93%  "

Set-TextValue $ws.Range("D15") "0.000009777"
Set-TextValue $ws.Range("E15") "  -6.35%  "

Set-TextValue $ws.Range("D16") "76.95"
Set-TextValue $ws.Range("E16") "  -5.38%  "

Set-TextValue $ws.Range("D17") "28.848.59"
Set-TextValue $ws.Range("E17") "  -1.58%  "

Set-TextValue $ws.Range("D18") "5.540"
Set-TextValue $ws.Range("E18") "  -11.17%  "

Set-TextValue $ws.Range("D19") "214.74"
Set-TextValue $ws.Range("E19") "  -8.94%  "

Set-TextValue $ws.Range("E20") "  +0.30%  "

Set-TextValue $ws.Range("D21") "11.48"
Set-TextValue $ws.Range("E21") "  -5.99%  "

Set-TextValue $ws.Range("D22") "6.861"
Set-TextValue $ws.Range("E22") "  -5.34%  "

Set-TextValue $ws.Range("D23") "1.005"
Set-TextValue $ws.Range("E23") "  +0.35%  "

Set-TextValue $ws.Range("D24") "156.35"

Set-TextValue $ws.Range("D25") "7.916"
Set-TextValue $ws.Range("E25") "  -5.86%  "

Set-TextValue $ws.Range("D26") "0.1281"
Set-TextValue $ws.Range("E26") "  -3.99%  "

Set-TextValue $ws.Range("D27") "16.42"
Set-TextValue $ws.Range("E27") "  -4.88%  "

Set-TextValue $ws.Range("D28") "0.06516"
Set-TextValue $ws.Range("E28") "  -8.50%  "

Set-TextValue $ws.Range("D29") "1.413"
Set-TextValue $ws.Range("E29") "  -4.46%  "

Set-TextValue $ws.Range("D30") "1.440"
Set-TextValue $ws.Range("E30") "  -2.66%  "

Set-TextValue $ws.Range("D31") "3.826"
Set-TextValue $ws.Range("E31") "  -4.77%  "

Set-TextValue $ws.Range("D32") "3.761"
Set-TextValue $ws.Range("E32") "  -6.45%  "

Set-TextValue $ws.Range("B33") "ARBITRUM"
Set-TextValue $ws.Range("C33") "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
Set-TextValue $ws.Range("D33") "1.090"
Set-TextValue $ws.Range("E33") "  -5.22%  "

Set-TextValue $ws.Range("B34") "LidoDAOToken"
Set-TextValue $ws.Range("C34") "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
Set-TextValue $ws.Range("D34") "1.721"
Set-TextValue $ws.Range("E34") "  -3.69%  "

Set-TextValue $ws.Range("D35") "0.6447"
Set-TextValue $ws.Range("E35") "  -6.34%  "

Set-TextValue $ws.Range("D36") "2.537"
Set-TextValue $ws.Range("E36") "  -1.72%  "

Set-TextValue $ws.Range("D37") "2.742"
Set-TextValue $ws.Range("E37") "  -1.55%  "

Set-TextValue $ws.Range("D38") "0.01752"
Set-TextValue $ws.Range("E38") "  -3.91%  "

Set-TextValue $ws.Range("D39") "6.494"
Set-TextValue $ws.Range("E39") "  -3.56%  "

Set-TextValue $ws.Range("D40") "1.137.00"
Set-TextValue $ws.Range("E40") "  -7.86%  "

Set-TextValue $ws.Range("D41") "0.8850"
Set-TextValue $ws.Range("E41") "  -6.27%  "

Set-TextValue $ws.Range("E42") "  +0.29%  "

Set-TextValue $ws.Range("D43") "1.988.54"
Set-TextValue $ws.Range("E43") "  -0.42%  "

Set-TextValue $ws.Range("D44") "99.89"
Set-TextValue $ws.Range("E44") "  -1.28%  "

Set-TextValue $ws.Range("D45") "61.53"
Set-TextValue $ws.Range("E45") "  -5.39%  "

Set-TextValue $ws.Range("E46") "  -3.45%  "

Set-TextValue $ws.Range("D47") "1.596"
Set-TextValue $ws.Range("E47") "  -5.07%  "

Set-TextValue $ws.Range("D48") "8.462"
Set-TextValue $ws.Range("E48") "  -4.78%  "

Set-TextValue $ws.Range("D49") "0.4541"
Set-TextValue $ws.Range("E49") "  -0.66%  "

Set-TextValue $ws.Range("D50") "0.05505"
Set-TextValue $ws.Range("E50") "  -2.45%  "

Set-TextValue $ws.Range("D51") "6.392"
Set-TextValue $ws.Range("E51") "  -7.58%  "
